# Updated TPM-based NATMI metrics for Sema3b-Nrp1 ligand-receptor pair
# Recomputed ligand/receptor expression values and their derived specificity
# and edge-weight metrics after refreshing the underlying TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "G"=0.6881063333333334; "H"=2.064319; "I"=0.04678220357266529; "J"=0.04678220357266529; "M"=123.2806423333333; "N"=369.841927; "O"=0.6241574062367528; "P"=0.6241574062367526; "Q"=84.83019076696812; "R"=763.4717169027131; "S"=0.02919945883995452; "T"=0.02919945883995451 }
    3 = @{ "G"=0.6881063333333334; "H"=2.064319; "I"=0.04678220357266529; "J"=0.04678220357266529; "O"=0.2392728888301323; "P"=0.2392728888301322; "Q"=32.51994545286944; "R"=292.679509075825; "S"=0.01119371299467096; "T"=0.01119371299467096 }
    4 = @{ "G"=0.6881063333333334; "H"=2.064319; "I"=0.04678220357266529; "J"=0.04678220357266529; "O"=0.136569704933115; "P"=0.136569704933115; "Q"=18.56139814524645; "R"=167.052583307218; "S"=0.006389031738039819; "T"=0.006389031738039818 }
    5 = @{ "I"=0.5086503334983149; "J"=0.5086503334983149; "M"=123.2806423333333; "N"=369.841927; "O"=0.6241574062367528; "P"=0.6241574062367526; "Q"=922.3358783714025; "R"=8301.022905342623; "S"=0.3174778728377675; "T"=0.3174778728377675 }
    6 = @{ "I"=0.5086503334983149; "J"=0.5086503334983149; "O"=0.2392728888301323; "P"=0.2392728888301322; "S"=0.121706234700552; "T"=0.121706234700552 }
    7 = @{ "I"=0.5086503334983149; "J"=0.5086503334983149; "O"=0.136569704933115; "P"=0.136569704933115; "S"=0.06946622595999542; "T"=0.06946622595999541 }
    8 = @{ "I"=0.4445674629290199; "J"=0.4445674629290199; "M"=123.2806423333333; "N"=369.841927; "O"=0.6241574062367528; "P"=0.6241574062367526; "Q"=806.1343803627761; "R"=7255.209423264985; "S"=0.2774800745590308; "T"=0.2774800745590307 }
    9 = @{ "I"=0.4445674629290199; "J"=0.4445674629290199; "O"=0.2392728888301323; "P"=0.2392728888301322; "S"=0.1063729411349093; "T"=0.1063729411349093 }
    10 = @{ "I"=0.4445674629290199; "J"=0.4445674629290199; "O"=0.136569704933115; "P"=0.136569704933115; "S"=0.06071444723507981; "T"=0.06071444723507979 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $updates[$row][$col]
    }
}
